# Insert a new price record as row 504 ("Hortaliza, Macroferia Regional de
# Talca - Papa" weekly data dump). This pushes the previous rows 504:550
# down to 505:551 and extends the used range to A1:R551.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(504).Insert()

$ws.Range("A504").Value = 5
$ws.Range("B504").Value = "Macroferia Regional de Talca"
$ws.Range("C504").Value = "Maule"
$ws.Range("D504").Value = 44769
$ws.Range("E504").Value = 7
$ws.Range("F504").Value = 100114001
$ws.Range("G504").Value = "Papa"
$ws.Range("H504").Value = "Patagonia"
$ws.Range("I504").Value = "1a (cosecha)"
$ws.Range("J504").Value = 1500
$ws.Range("K504").Value = 7000
$ws.Range("L504").Value = 7000
$ws.Range("M504").Value = 7000
$ws.Range("N504").Value = "`$/saco 25 kilos"
$ws.Range("O504").Value = "Región de La Araucanía"
$ws.Range("P504").Value = 280
$ws.Range("Q504").Value = 25
$ws.Range("R504").Value = "Hortaliza"
